$d = $word.ActiveDocument

# --- Locate the "C.C: Documento de identidad" paragraph in the table (the
#     one whose run immediately after "C.C: " is "Documento de identidad")
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("C.C: ")
$ccStart = $rng.Start

$rng2 = $d.Range($rng.End, $d.Content.End)
$found2 = $rng2.Find.Execute("Documento de identidad")
$docEnd = $rng2.End

$full = $d.Range($ccStart, $docEnd)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1C1F6368" w14:textId="2B926550" w:rsidR="00F7036C" w:rsidRDefault="003844FE">
<w:pPr>
<w:ind w:left="0" w:hanging="2"/>
<w:jc w:val="both"/>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
</w:rPr>
<w:t>C.C:</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r w:rsidR="00C8324B" w:rsidRPr="00C8324B">
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
<w:b/>
<w:bCs/>
<w:color w:val="000000" w:themeColor="text1"/>
</w:rPr>
<w:t>Documento_trabajador</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$full.InsertXML($xmlFrag)
Write-Output "C.C block rewritten"
